# Fruta / hortaliza, semanal
# Re-order the weekly price records in rows 3-7 (dates, volumes, prices) to
# reflect the updated weekly sequence, leaving rows 2 and 8 untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44253
$ws.Range("M3").Value = 90
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 13000
$ws.Range("P3").Value = 12667
$ws.Range("S3").Value = 905

# Row 4
$ws.Range("D4").Value = 44210
$ws.Range("M4").Value = 70
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 11000
$ws.Range("P4").Value = 10357
$ws.Range("S4").Value = 740

# Row 5
$ws.Range("D5").Value = 44216
$ws.Range("M5").Value = 55
$ws.Range("N5").Value = 11000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 11545
$ws.Range("S5").Value = 825

# Row 6
$ws.Range("D6").Value = 44232
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 11000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 11583
$ws.Range("S6").Value = 827

# Row 7
$ws.Range("D7").Value = 44229
$ws.Range("M7").Value = 55
$ws.Range("N7").Value = 11000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 11364
$ws.Range("S7").Value = 812
